# "final code of the framework"
# Populate Sheet1 with test-account credentials (email + password), turn the
# email column into mailto: hyperlinks, and style the password column with a
# monospaced font so it's easy to read/copy during manual QA.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data -----------------------------------------------------------------
$ws.Range("A1").Value = "saradvd25@hotmail.com"
$ws.Range("B1").Value = "Pa55word"
$ws.Range("A2").Value = "owenck25@gmail.com"
$ws.Range("B2").Value = "Pa55word"
$ws.Range("A3").Value = "darasl50@gmail.com"
$ws.Range("B3").Value = "Pa55word"

# --- mailto hyperlinks on the email column ---------------------------------
$ws.Hyperlinks.Add($ws.Range("A1"), "mailto:saradvd25@hotmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:owenck25@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:darasl50@gmail.com") | Out-Null

# --- password column font (Consolas 12, blue-ish) --------------------------
$pwdRange = $ws.Range("B1:B3")
$pwdRange.Font.Name = "Consolas"
$pwdRange.Font.Size = 12
$pwdRange.Font.Color = 16711722

# --- column widths / row heights to fit the content -------------------------
$ws.Columns.Item(1).ColumnWidth = 30.6
$ws.Columns.Item(2).ColumnWidth = 17

$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(3).RowHeight = 15.75

# --- leave the selection on the last row, like the authored workbook -------
$ws.Range("A3").Select() | Out-Null
